# Auto-generated edit script: updates market-price-derived columns (H-N)
# for specific leve rows across multiple sheets, per the scheduled runner's
# refreshed Universalis price data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 8384280
$ws.Range("I70").Value = 33534080
$ws.Range("J70").Value = 1012.8
$ws.Range("K70").Value = 100602240
$ws.Range("L70").Value = 3038.4
$ws.Range("M70").Value = -100601970
$ws.Range("N70").Value = -3578.4
# Row 73
$ws.Range("H73").Value = 8384280
$ws.Range("I73").Value = 33534080
$ws.Range("J73").Value = 1012.8
$ws.Range("K73").Value = 100602240
$ws.Range("L73").Value = 3038.4
$ws.Range("M73").Value = -100601304
$ws.Range("N73").Value = -4910.4
# Row 112
$ws.Range("H112").Value = 6337.7144
$ws.Range("J112").Value = 6994.24
$ws.Range("L112").Value = 20982.72
$ws.Range("N112").Value = -23198.72
# Row 138
$ws.Range("H138").Value = 2444895.2
$ws.Range("I138").Value = 6254161
$ws.Range("J138").Value = 6965.16
$ws.Range("K138").Value = 18762483
$ws.Range("L138").Value = 20895.48
$ws.Range("M138").Value = -18757343
$ws.Range("N138").Value = -31175.48

$ws = $wb.Worksheets.Item("ARM")
# Row 44
$ws.Range("H44").Value = 37949
$ws.Range("J44").Value = 37949
$ws.Range("L44").Value = 37949
$ws.Range("N44").Value = -38925
# Row 45
$ws.Range("H45").Value = 15244.286
$ws.Range("I45").Value = 555
$ws.Range("J45").Value = 21120
$ws.Range("K45").Value = 555
$ws.Range("L45").Value = 21120
$ws.Range("M45").Value = -178
$ws.Range("N45").Value = -21874
# Row 102
$ws.Range("H102").Value = 126751.25
$ws.Range("I102").Value = 1842
$ws.Range("J102").Value = 334933.34
$ws.Range("K102").Value = 1842
$ws.Range("L102").Value = 334933.34
$ws.Range("M102").Value = -220
$ws.Range("N102").Value = -338177.34
# Row 132
$ws.Range("H132").Value = 2566.7778
$ws.Range("I132").Value = 2234.1853
$ws.Range("K132").Value = 6702.5559
$ws.Range("M132").Value = -4172.5559

$ws = $wb.Worksheets.Item("BSM")
# Row 57
$ws.Range("H57").Value = 52239.5
$ws.Range("I57").Value = 5709
$ws.Range("J57").Value = 98770
$ws.Range("K57").Value = 5709
$ws.Range("L57").Value = 98770
$ws.Range("M57").Value = -4989
$ws.Range("N57").Value = -100210
# Row 99
$ws.Range("H99").Value = 1779.5333
$ws.Range("I99").Value = 1188.1111
$ws.Range("J99").Value = 2666.6667
$ws.Range("K99").Value = 1188.1111
$ws.Range("L99").Value = 2666.6667
$ws.Range("M99").Value = 309.8888999999999
$ws.Range("N99").Value = -5662.6667
# Row 105
$ws.Range("H105").Value = 799999.6
$ws.Range("I105").Value = 837961.5
$ws.Range("J105").Value = 2800
$ws.Range("K105").Value = 837961.5
$ws.Range("L105").Value = 2800
$ws.Range("M105").Value = -836214.5
$ws.Range("N105").Value = -6294
# Row 136
$ws.Range("H136").Value = 52239.5
$ws.Range("I136").Value = 5709
$ws.Range("J136").Value = 98770
$ws.Range("K136").Value = 5709
$ws.Range("L136").Value = 98770
$ws.Range("M136").Value = -609
$ws.Range("N136").Value = -108970

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1276.8918
$ws.Range("I58").Value = 1342.1786
$ws.Range("J58").Value = 1073.7778
$ws.Range("K58").Value = 1342.1786
$ws.Range("L58").Value = 1073.7778
$ws.Range("M58").Value = -1139.1786
$ws.Range("N58").Value = -1479.7778
# Row 136
$ws.Range("H136").Value = 1276.8918
$ws.Range("I136").Value = 1342.1786
$ws.Range("J136").Value = 1073.7778
$ws.Range("K136").Value = 4026.5358
$ws.Range("L136").Value = 3221.3334
$ws.Range("M136").Value = -1476.5358
$ws.Range("N136").Value = -8321.3334

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 3634.8538
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 3703.225
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 11109.675
$ws.Range("M113").Value = -530
$ws.Range("N113").Value = -15449.675
# Row 124
$ws.Range("H124").Value = 2414.1177
$ws.Range("I124").Value = 782.2222
$ws.Range("J124").Value = 4250
$ws.Range("K124").Value = 2346.6666
$ws.Range("L124").Value = 12750
$ws.Range("M124").Value = 2563.3334
$ws.Range("N124").Value = -22570
# Row 131
$ws.Range("H131").Value = 18185130
$ws.Range("J131").Value = 18869460
$ws.Range("L131").Value = 56608380
$ws.Range("N131").Value = -56618460
# Row 132
$ws.Range("H132").Value = 1984.5
$ws.Range("I132").Value = 1548.2858
$ws.Range("J132").Value = 2106.64
$ws.Range("K132").Value = 13934.5722
$ws.Range("L132").Value = 18959.76
$ws.Range("M132").Value = -11404.5722
$ws.Range("N132").Value = -24019.76
# Row 136
$ws.Range("H136").Value = 5396.316
$ws.Range("I136").Value = 1424
$ws.Range("J136").Value = 6815
$ws.Range("K136").Value = 4272
$ws.Range("L136").Value = 20445
$ws.Range("M136").Value = 828
$ws.Range("N136").Value = -30645

$ws = $wb.Worksheets.Item("GSM")
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1143.875
$ws.Range("I22").Value = 150
$ws.Range("J22").Value = 2137.75
$ws.Range("K22").Value = 150
$ws.Range("L22").Value = 2137.75
$ws.Range("M22").Value = 145
$ws.Range("N22").Value = -2727.75
# Row 27
$ws.Range("H27").Value = 1143.875
$ws.Range("I27").Value = 150
$ws.Range("J27").Value = 2137.75
$ws.Range("K27").Value = 150
$ws.Range("L27").Value = 2137.75
$ws.Range("M27").Value = -43
$ws.Range("N27").Value = -2351.75
# Row 46
$ws.Range("H46").Value = 1277.8
$ws.Range("I46").Value = 1057.6
$ws.Range("J46").Value = 1387.9
$ws.Range("K46").Value = 1057.6
$ws.Range("L46").Value = 1387.9
$ws.Range("M46").Value = -869.5999999999999
$ws.Range("N46").Value = -1763.9
# Row 55
$ws.Range("H55").Value = 338.41666
$ws.Range("I55").Value = 360.0909
$ws.Range("K55").Value = 360.0909
$ws.Range("M55").Value = -187.0909

$ws = $wb.Worksheets.Item("WVR")
# Row 112
$ws.Range("H112").Value = 127938.7
$ws.Range("J112").Value = 127938.7
$ws.Range("L112").Value = 127938.7
$ws.Range("N112").Value = -130892.7
# Row 123
$ws.Range("H123").Value = 54429
$ws.Range("J123").Value = 54429
$ws.Range("L123").Value = 54429
$ws.Range("N123").Value = -64229
